$d = $word.ActiveDocument

# --- Paragraph "Estimado Sr. Antonio..." -------------------------------
# The paragraph currently reads (as separate runs, split by proofErr
# markers from the grammar checker):
#   <w:proofErr gramStart/> "Estimado" <w:proofErr gramEnd/> " Sr. Antonio"
# We want it merged into a single run "Estimado Sr. Antonio" with both
# proofErr markers gone. A find/replace that starts exactly at the
# paragraph's first character leaves a proofErr sitting right at the
# match boundary untouched, so we first insert a temporary marker run
# in front of it (making the leading proofErr sit *between* two runs,
# i.e. strictly inside the next match) and fold it back out again while
# replacing.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text -like "Estimado Sr. Antonio*") {
        $target = $cand
        break
    }
}
if ($target -eq $null) { $target = $d.Paragraphs(3) }
$target.Range.InsertBefore("@@MARK@@")
$d.Content.Find.Execute("@@MARK@@Estimado Sr. Antonio", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Estimado Sr. Antonio", 2)

# --- Paragraph "Fecha del evento: ..." ----------------------------------
# Same underlying pattern (text split by proofErr markers around
# "Sábado 18 de Diciembre"); here the proofErr pair sits entirely inside
# the run span we are replacing, so a plain find/replace already folds
# it away cleanly.
$d.Content.Find.Execute("Fecha del evento: Sábado 18 de Diciembre de 2023", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "Fecha del evento: Sábado 18 de Diciembre de 2023", 2)
